# Timer funzioni sincrone OK
# Append two empty paragraphs (same style/spacing as the surrounding body
# text) at the very end of the document, after the last existing
# paragraph ("Capire se è proprio necessario...") and before the
# section properties.

$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
